$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20..126 down to 21..127
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44473
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112021
$ws.Range("G20").Value = "Ají"
$ws.Range("H20").Value = "Americana (o)"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 85000
$ws.Range("L20").Value = 85000
$ws.Range("M20").Value = 85000
$ws.Range("N20").Value = "$/caja 25 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 3400
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
